$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H3").Value = 27171.428
$ws.Range("J3").Value = 27171.428
$ws.Range("L3").Value = 27171.428
$ws.Range("N3").Value = -27399.428
$ws.Range("H51").Value = 4623.5454
$ws.Range("I51").Value = 7737.25
$ws.Range("J51").Value = 2844.2856
$ws.Range("K51").Value = 7737.25
$ws.Range("L51").Value = 2844.2856
$ws.Range("M51").Value = -7253.25
$ws.Range("N51").Value = -3812.2856
$ws.Range("H100").Value = 90910780
$ws.Range("I100").Value = 142857970
$ws.Range("K100").Value = 142857970
$ws.Range("M100").Value = -142857429
$ws.Range("H102").Value = 27171.428
$ws.Range("J102").Value = 27171.428
$ws.Range("L102").Value = 27171.428
$ws.Range("N102").Value = -33661.428
$ws.Range("H132").Value = 3098.9285
$ws.Range("I132").Value = 3362.8
$ws.Range("J132").Value = 900
$ws.Range("K132").Value = 10088.4
$ws.Range("L132").Value = 2700
$ws.Range("M132").Value = -7558.400000000001
$ws.Range("N132").Value = -7760
$ws.Range("H135").Value = 17862720
$ws.Range("I135").Value = 550.26086
$ws.Range("K135").Value = 4952.34774
$ws.Range("M135").Value = -2417.34774
$ws.Range("H137").Value = 1999.6
$ws.Range("I137").Value = 1590.8182
$ws.Range("K137").Value = 4772.4546
$ws.Range("M137").Value = -2222.4546
$ws.Range("H138").Value = 2165.2346
$ws.Range("I138").Value = 2019.1333
$ws.Range("J138").Value = 2198.4395
$ws.Range("K138").Value = 6057.3999
$ws.Range("L138").Value = 6595.318499999999
$ws.Range("M138").Value = -917.3999000000003
$ws.Range("N138").Value = -16875.3185
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 5050.0166
$ws.Range("I32").Value = 5129.608
$ws.Range("J32").Value = 4644.1
$ws.Range("K32").Value = 5129.608
$ws.Range("L32").Value = 4644.1
$ws.Range("M32").Value = -4842.608
$ws.Range("N32").Value = -5218.1
$ws.Range("H45").Value = 2303.0244
$ws.Range("I45").Value = 2187.7856
$ws.Range("J45").Value = 2551.2307
$ws.Range("K45").Value = 2187.7856
$ws.Range("L45").Value = 2551.2307
$ws.Range("M45").Value = -1810.7856
$ws.Range("N45").Value = -3305.2307
$ws.Range("H74").Value = 200000800
$ws.Range("I74").Value = 250000770
$ws.Range("J74").Value = 890
$ws.Range("K74").Value = 250000770
$ws.Range("L74").Value = 890
$ws.Range("M74").Value = -249999896
$ws.Range("N74").Value = -2638
$ws.Range("H77").Value = 200000800
$ws.Range("I77").Value = 250000770
$ws.Range("J77").Value = 890
$ws.Range("K77").Value = 1250003850
$ws.Range("L77").Value = 4450
$ws.Range("M77").Value = -1249999482
$ws.Range("N77").Value = -13186
$ws.Range("H102").Value = 1750
$ws.Range("I102").Value = 0
$ws.Range("K102").Value = 0
$ws.Range("H122").Value = 2799
$ws.Range("I122").Value = 2237.4
$ws.Range("J122").Value = 4671
$ws.Range("K122").Value = 6712.200000000001
$ws.Range("L122").Value = 14013
$ws.Range("M122").Value = -4262.200000000001
$ws.Range("N122").Value = -18913
$ws.Range("H132").Value = 31845.9
$ws.Range("I132").Value = 1662.1351
$ws.Range("J132").Value = 117753.54
$ws.Range("K132").Value = 4986.4053
$ws.Range("L132").Value = 353260.62
$ws.Range("M132").Value = -2456.4053
$ws.Range("N132").Value = -358320.62
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 4751.8887
$ws.Range("I105").Value = 6692
$ws.Range("K105").Value = 6692
$ws.Range("M105").Value = -4945
$ws.Range("H134").Value = 8232.615
$ws.Range("I134").Value = 9002.182000000001
$ws.Range("K134").Value = 27006.546
$ws.Range("M134").Value = -24471.546
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 25224.428
$ws.Range("I58").Value = 1520
$ws.Range("J58").Value = 84485.5
$ws.Range("K58").Value = 1520
$ws.Range("L58").Value = 84485.5
$ws.Range("M58").Value = -1317
$ws.Range("N58").Value = -84891.5
$ws.Range("H76").Value = 0
$ws.Range("I76").Value = 0
$ws.Range("K76").Value = 0
$ws.Range("H79").Value = 0
$ws.Range("I79").Value = 0
$ws.Range("K79").Value = 0
$ws.Range("H105").Value = 12500795
$ws.Range("I105").Value = 12500795
$ws.Range("K105").Value = 12500795
$ws.Range("M105").Value = -12499048
$ws.Range("H132").Value = 24060.25
$ws.Range("I132").Value = 30006.666
$ws.Range("J132").Value = 6221
$ws.Range("K132").Value = 90019.99800000001
$ws.Range("L132").Value = 18663
$ws.Range("M132").Value = -87489.99800000001
$ws.Range("N132").Value = -23723
$ws.Range("H136").Value = 25224.428
$ws.Range("I136").Value = 1520
$ws.Range("J136").Value = 84485.5
$ws.Range("K136").Value = 4560
$ws.Range("L136").Value = 253456.5
$ws.Range("M136").Value = -2010
$ws.Range("N136").Value = -258556.5
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 5000913
$ws.Range("I4").Value = 1369.75
$ws.Range("K4").Value = 4109.25
$ws.Range("M4").Value = -3997.25
$ws.Range("H131").Value = 819.33
$ws.Range("I131").Value = 800
$ws.Range("J131").Value = 819.52527
$ws.Range("K131").Value = 2400
$ws.Range("L131").Value = 2458.57581
$ws.Range("M131").Value = 2640
$ws.Range("N131").Value = -12538.57581
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H5").Value = 0
$ws.Range("J5").Value = 0
$ws.Range("N5").Value = 0
$ws.Range("H70").Value = 23600
$ws.Range("I70").Value = 50000
$ws.Range("J70").Value = 6000
$ws.Range("K70").Value = 50000
$ws.Range("L70").Value = 6000
$ws.Range("M70").Value = -49730
$ws.Range("N70").Value = -6540
$ws.Range("H73").Value = 23600
$ws.Range("I73").Value = 50000
$ws.Range("J73").Value = 6000
$ws.Range("K73").Value = 50000
$ws.Range("L73").Value = 6000
$ws.Range("M73").Value = -49064
$ws.Range("N73").Value = -7872
$ws.Range("H80").Value = 3927.7273
$ws.Range("I80").Value = 3601.25
$ws.Range("J80").Value = 4114.2856
$ws.Range("K80").Value = 3601.25
$ws.Range("L80").Value = 4114.2856
$ws.Range("M80").Value = -2603.25
$ws.Range("N80").Value = -6110.2856
$ws.Range("H83").Value = 3927.7273
$ws.Range("I83").Value = 3601.25
$ws.Range("J83").Value = 4114.2856
$ws.Range("K83").Value = 18006.25
$ws.Range("L83").Value = 20571.428
$ws.Range("M83").Value = -13014.25
$ws.Range("N83").Value = -30555.428
$ws.Range("H122").Value = 266667460
$ws.Range("I122").Value = 83334320
$ws.Range("J122").Value = 1000000000
$ws.Range("K122").Value = 250002960
$ws.Range("L122").Value = 3000000000
$ws.Range("M122").Value = -250000510
$ws.Range("N122").Value = -3000004900
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 756620.5
$ws.Range("I122").Value = 1309388.9
$ws.Range("K122").Value = 3928166.7
$ws.Range("M122").Value = -3925716.7
$ws.Range("H132").Value = 2106
$ws.Range("I132").Value = 1828.2
$ws.Range("J132").Value = 2800.5
$ws.Range("K132").Value = 5484.6
$ws.Range("L132").Value = 8401.5
$ws.Range("M132").Value = -2954.6
$ws.Range("N132").Value = -13461.5
$ws.Range("H136").Value = 29491.555
$ws.Range("I136").Value = 39672.152
$ws.Range("J136").Value = 3022
$ws.Range("K136").Value = 119016.456
$ws.Range("L136").Value = 9066
$ws.Range("M136").Value = -116466.456
$ws.Range("N136").Value = -14166
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H2").Value = 19902
$ws.Range("J2").Value = 0
$ws.Range("L2").Value = 0
$ws.Range("H56").Value = 2642.5
$ws.Range("I56").Value = 2642.5
$ws.Range("K56").Value = 2642.5
$ws.Range("M56").Value = -1928.5
$ws.Range("H122").Value = 1318.2307
$ws.Range("I122").Value = 1188
$ws.Range("J122").Value = 1399.625
$ws.Range("K122").Value = 3564
$ws.Range("L122").Value = 4198.875
$ws.Range("M122").Value = -1114
$ws.Range("N122").Value = -9098.875
$ws.Range("H133").Value = 43476.668
$ws.Range("J133").Value = 43476.668
$ws.Range("L133").Value = 43476.668
$ws.Range("N133").Value = -53596.668
$ws.Range("H140").Value = 45672.5
$ws.Range("J140").Value = 45672.5
$ws.Range("L140").Value = 45672.5
$ws.Range("N140").Value = -56032.5
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("M102").ClearContents()
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("M76").ClearContents()
$ws.Range("M79").ClearContents()
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("N5").ClearContents()
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("N2").ClearContents()
